$wb = $excel.ActiveWorkbook

$metadata = $wb.Worksheets.Item("Metadata")
$concepts = $wb.Worksheets.Item("Concepts")

# Update the Date value (row 8, column B) on the Metadata sheet
$metadata.Range("B8").Value = "2021-12-08T13:27:26-05:00"

# Replace the "todo" Description (row 13, column B) on the Metadata sheet
$metadata.Range("B13").Value = "Categories for SPLASCH observations"

# Replace "todo definition" entries on the Concepts sheet with actual definitions
$concepts.Range("D2").Value = "Category code for spoken language comprehension observation"
$concepts.Range("D3").Value = "Category code for spoken language expression observation"
$concepts.Range("D4").Value = "Category code for swallowing observation"
